$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the country data rows (COVID-19 case counts refresh) ---

# Alemania (row 8)
$ws.Range("B8").Value = 148746
$ws.Range("C8").Value = 293
$ws.Range("E8").Value = 44244
$ws.Range("G8").Value = 16
$ws.Range("H8").Value = 5102

# Portugal (row 19)
$ws.Range("B19").Value = 21982
$ws.Range("C19").Value = 603
$ws.Range("D19").Value = 1143
$ws.Range("E19").Value = 20054
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = 785

# India (row 20)
$ws.Range("B20").Value = 20471
$ws.Range("C20").Value = 391
$ws.Range("E20").Value = 15843
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 652

# Arabia Saudita (row 26)
$ws.Range("B26").Value = 12772
$ws.Range("C26").Value = 1141
$ws.Range("D26").Value = 1812
$ws.Range("E26").Value = 10846
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 114

# --- Eslovaquia / Republica de Macedonia swap positions (row 78/79) ---
# Row 78 keeps "Eslovaquia" label -> becomes "Republica de Macedonia" with fresh data
# Row 79 keeps "Republica de Macedonia" label -> becomes "Eslovaquia" with the data that
# previously belonged to row 78.
$ws.Range("A78").Value = "Republica de Macedonia"
$ws.Range("B78").Value = 1259
$ws.Range("C78").Value = 28
$ws.Range("D78").Value = 272
$ws.Range("E78").Value = 931
$ws.Range("F78").Value = 8
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 56

$ws.Range("A79").Value = "Eslovaquia"
$ws.Range("B79").Value = 1244
$ws.Range("C79").Value = 45
$ws.Range("D79").Value = 284
$ws.Range("E79").Value = 946
$ws.Range("F79").Value = 7
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 14

# --- Update "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 14:52"
